$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 and 45: coin order swap (Bittensor <-> InjectiveProtocol) with updated values
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.62%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "351.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.93%  "

$ws.Range("D2").Value = "67.952.40"
$ws.Range("E2").Value = "  +3.39%  "

$ws.Range("D3").Value = "3.282.46"
$ws.Range("E3").Value = "  +3.50%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.12%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.02%  "

$ws.Range("D9").Value = "3.283.83"
$ws.Range("E9").Value = "  +3.57%  "

$ws.Range("E10").Value = "  +7.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.87%  "

$ws.Range("E12").Value = "  +6.50%  "

$ws.Range("D13").Value = "3.851.23"
$ws.Range("E13").Value = "  +3.59%  "

$ws.Range("E14").Value = "  +1.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.88%  "

$ws.Range("D16").Value = "67.944.21"
$ws.Range("E16").Value = "  +3.39%  "

$ws.Range("E17").Value = "  +3.40%  "

$ws.Range("D18").Value = "3.278.25"
$ws.Range("E18").Value = "  +3.25%  "

$ws.Range("E19").Value = "  +1.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.49%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.515"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.83%  "

$ws.Range("E26").Value = "  +5.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.46%  "

$ws.Range("E28").Value = "  +2.39%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").Value = "  +3.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.02%  "

$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("E35").Value = "  +4.61%  "

$ws.Range("E36").Value = "  +5.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.04%  "

$ws.Range("E38").Value = "  +1.95%  "

$ws.Range("E39").Value = "  +3.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.83%  "

$ws.Range("E41").Value = "  +9.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.55%  "

$ws.Range("D46").Value = "2.676.72"
$ws.Range("E46").Value = "  +0.73%  "

$ws.Range("E47").Value = "  +2.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0683"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0285"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.34%  "

$ws.Range("E50").Value = "  +5.21%  "
